$d = $word.ActiveDocument

function Replace-ExactText($doc, $oldText, $newText) {
    $rng = $doc.Content.Duplicate
    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text to replace was not found: $oldText"
    }
    # Assigning straight to Range.Text (rather than doing the substitution
    # through Find.Execute's Replace argument) avoids Word's "smart quotes"
    # AutoFormat kicking in and mangling the literal " characters that are
    # part of this JSON payload.
    $rng.Text = $newText
}

# 1) Update the three JSON field values that changed.
Replace-ExactText $d "這是他創意大腦最重要的產物。" "其最終目的是完全掌握心靈對物質世界的控制，將自然的力量運用於人類的需求。"

Replace-ExactText $d "It is the most important product of his creative brain." "Its ultimate purpose is the complete mastery of mind over the material world, the harnessing of the forces of nature to human needs."

Replace-ExactText $d "The translation is grammatically correct and conveys the meaning accurately. However, the word choice for 'creative brain' is too literal.  A more natural translation would be '創意' instead of '創意大腦'." 'Translation is too formal and wordy. \"掌握心靈對物質世界的控制\" is awkward. \"Harnessing\" should be translated as \"利用\" or \"駕馭\".'

# 2) The seven pieces of JSON text currently live in a single paragraph,
#    separated by manual line breaks (Chr(11), i.e. <w:br/>). Turn each
#    manual line break into a real paragraph mark (Chr(13)) so every piece
#    becomes its own paragraph.
$d.Content.Find.Execute(
    [string][char]11,
    $false, $false, $false, $false, $false, $true, 1, $false,
    [string][char]13,
    2)
